$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 517; this shifts existing rows 517..588
# down to 518..589 (matching the diff's row-shift pattern) and extends
# the used range/dimension automatically.
$ws.Rows.Item(517).Insert()

# Populate the newly inserted row 517 with a weekly record for
# Brocoli / Macroferia Regional de Talca, mirroring the other rows'
# static columns and the new price data from the commit.
$ws.Range("A517").Value = 5
$ws.Range("B517").Value = "Macroferia Regional de Talca"
$ws.Range("C517").Value = "Maule"
$ws.Range("D517").Value = 45124
$ws.Range("E517").Value = 7
$ws.Range("F517").Value = 100112023
$ws.Range("G517").Value = "Brócoli"
$ws.Range("H517").Value = "Sin especificar"
$ws.Range("I517").Value = "Primera"
$ws.Range("J517").Value = 3000
$ws.Range("K517").Value = 700
$ws.Range("L517").Value = 700
$ws.Range("M517").Value = 700
$ws.Range("N517").Value = "$/unidad"
$ws.Range("O517").Value = "Región del Maule"
$ws.Range("P517").Value = 700
$ws.Range("Q517").Value = 1
$ws.Range("R517").Value = "Hortaliza"
